$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 6 (C6:I6) values to reflect the redistributed single-EV power values
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = 0.95
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 25
$ws.Range("I6").Value = 0.9

# Update the selection to match the edited range
$ws.Activate()
$ws.Range("C6:I6").Select()
